$d = $word.ActiveDocument

# --- 1. Table width: auto -> 100% (pct, w:w=5000) ---
$t = $d.Tables.Item(1)
$t.PreferredWidthType = 2   # wdPreferredWidthPercent
$t.PreferredWidth = 250     # PreferredWidth(250) * 20 = 5000 => 100%

# --- 2. Fill in the rubric row (row 2) with left-justified descriptions ---
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$cell1Xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Compact`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Response directly addresses the prompt with specific details from the readings and NotebookLM. Includes concrete examples from field experience or teaching practice.</w:t></w:r></w:p>"
$t.Rows.Item(2).Cells.Item(1).Range.Paragraphs.Item(1).Range.InsertXML($cell1Xml)

$cell2Xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Compact`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Response addresses the prompt but lacks specific details or examples. May be vague or general.</w:t></w:r></w:p>"
$t.Rows.Item(2).Cells.Item(2).Range.Paragraphs.Item(1).Range.InsertXML($cell2Xml)

$cell3Xml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Compact`"/><w:jc w:val=`"left`"/></w:pPr><w:r><w:t xml:space=`"preserve`">No response or response does not address the prompt.</w:t></w:r></w:p>"
$t.Rows.Item(2).Cells.Item(3).Range.Paragraphs.Item(1).Range.InsertXML($cell3Xml)

Write-Output "Rubric table updated."
